$wb = $excel.ActiveWorkbook

# --- Step 1: insert a new "2022-Q1" sheet right before the "总计" sheet.
#     The "总计" sheet is removed and re-added after the new sheet so the
#     internal sheetId sequence matches what Excel itself would produce
#     (ids are handed out in creation order, not by position): the fresh
#     "2022-Q1" sheet becomes sheetId 6 and "总计" becomes sheetId 7.
$q4Sheet    = $wb.Worksheets.Item(5)
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totalName  = $totalSheet.Name
$totalSheet.Delete()

$q1Sheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1Sheet)
$totalSheet.Name = $totalName

# A reliable style donor: row 2 of the "2021-Q4" sheet carries the same
# header-row / index-column formatting ("s=2": bold, centered, thin box
# border) that every quarter sheet in this workbook uses.
$styleRow = $q4Sheet.Range("A2:H2")

function Copy-HeaderStyle($destRange) {
    $styleRow.Cells.Item(1, 1).Copy() | Out-Null
    $destRange.PasteSpecial(-4122) | Out-Null
}

# --- Step 2: populate the new "2022-Q1" fund-holding detail sheet.
$ws = $q1Sheet

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "1")
    Copy-HeaderStyle $cell
    $cell.Value = $headers[$i]
}

$data = @(
    @("400003", "东方精选混合",                 "10.47", "82.20", "3.51", "0.3675", 7),
    @("400001", "东方龙混合",                   "2.80",  "84.04", "3.90", "0.1092", 7),
    @("009937", "东方欣益一年持有期偏债混合A",   "3.39",  "31.18", "1.77", "0.0600", 7),
    @("009169", "湘财长兴灵活配置混合A",         "1.16",  "85.40", "2.94", "0.0341", 8),
    @("009170", "湘财长兴灵活配置混合C",         "0.46",  "85.40", "2.94", "0.0135", 8),
    @("009938", "东方欣益一年持有期偏债混合C",   "0.52",  "31.18", "1.77", "0.0092", 7)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $aCell = $ws.Range("A$row")
    Copy-HeaderStyle $aCell
    $aCell.Value = $r

    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $data[$r][0]
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = $data[$r][1]
    $ws.Range("C$row").Style = "Normal"

    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $data[$r][2]
    $ws.Range("D$row").Style = "Normal"

    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = $data[$r][3]
    $ws.Range("E$row").Style = "Normal"

    $ws.Range("F$row").NumberFormat = "@"
    $ws.Range("F$row").Value = $data[$r][4]
    $ws.Range("F$row").Style = "Normal"

    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = $data[$r][5]
    $ws.Range("G$row").Style = "Normal"

    $ws.Range("H$row").Value = $data[$r][6]
}

# --- Step 3: re-populate "总计" with a new leading "2022-Q1" row, shifting
#     the previously existing rows down by one.
$ts = $totalSheet

$totHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
$totCols = @("B", "C", "D")
for ($i = 0; $i -lt $totHeaders.Length; $i++) {
    $cell = $ts.Range($totCols[$i] + "1")
    Copy-HeaderStyle $cell
    $cell.Value = $totHeaders[$i]
}

$totData = @(
    @("2022-Q1", 6,  0.59),
    @("2021-Q4", 8,  0.32),
    @("2021-Q3", 4,  0.59),
    @("2021-Q2", 13, 2.19),
    @("2021-Q1", 13, 0.97),
    @("2020-Q4", 10, 4.25)
)

for ($r = 0; $r -lt $totData.Length; $r++) {
    $row = $r + 2
    $aCell = $ts.Range("A$row")
    Copy-HeaderStyle $aCell
    $aCell.Value = $r

    $ts.Range("B$row").Value = $totData[$r][0]
    $ts.Range("C$row").Value = $totData[$r][1]
    $ts.Range("D$row").Value = $totData[$r][2]
}
